# Update "想去人数" (F column) figures and one cover image URL (I column)
# on the "展览" sheet (rows 7-18) and the "全部类型" sheet (rows 8-19, same
# events shifted down by one row), matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $rowOffset = 0
    } else {
        $rowOffset = 1
    }

    # F column: 想去人数 updates
    $ws.Range("F" + (7 + $rowOffset)).Value = 561
    $ws.Range("F" + (8 + $rowOffset)).Value = 7907
    $ws.Range("F" + (9 + $rowOffset)).Value = 751
    $ws.Range("F" + (10 + $rowOffset)).Value = 222
    $ws.Range("F" + (11 + $rowOffset)).Value = 1095
    $ws.Range("F" + (12 + $rowOffset)).Value = 749
    $ws.Range("F" + (13 + $rowOffset)).Value = 28
    $ws.Range("F" + (15 + $rowOffset)).Value = 197
    $ws.Range("F" + (16 + $rowOffset)).Value = 40
    $ws.Range("F" + (18 + $rowOffset)).Value = 817

    # I column: cover image URL refresh for 安徽·赛马娘Only 2.0
    $ws.Range("I" + (16 + $rowOffset)).Value = "//i1.hdslb.com/bfs/openplatform/202405/ibcY9Edj1715235810905.jpeg"
}
